$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D): force text storage so numeric-looking strings
# (e.g. "1.001") are not auto-coerced to numbers by Excel, matching
# the original inline-string cell type. Resetting the style back to
# "Normal" afterwards avoids leaving a stray Text number-format on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.294.91"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.01%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.921.89"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.66%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.8143"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.80%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "244.39"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.12%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3257"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.14%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "27.27"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.78%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07238"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.69%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7945"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +7.19%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08117"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.45%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.934.38"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.33%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.441"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.83%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "94.42"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.61%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.300.94"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.02%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.27"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.22%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.088"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.76%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "250.09"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.69%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007865"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.54%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.186.48"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.28%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.212"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +20.25%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.11%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.002"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.08%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1664"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +19.21%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.516"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.92%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "168.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.00%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.06"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.79%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.157"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.32%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.372"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.51%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.553"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.72%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.355"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.05%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05752"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.78%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.153"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.67%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.305"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.78%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7493"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.37%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9988"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.07%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.729"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.31%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01962"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.06%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.819"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.27%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4513"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.31%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "74.79"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.50%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.987"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.09%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8571"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.78%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.931"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.87%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.041.45"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.49%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.000"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.02%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "103.37"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.83%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.123"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +11.34%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.666"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.82%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.946"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.86%  "

